$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plants")

# Set selected variables to not being tuned ("no" instead of "yes")
$ws.Range("E6").Value = "no"
$ws.Range("E7").Value = "no"
$ws.Range("E19").Value = "no"
$ws.Range("E20").Value = "no"

# Update the active selection/cell to D5
$ws.Range("D5").Select()
